$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D1) and "is_enabled_lbl" (E1) columns' header
# cells by deleting them and shifting the remaining headers (order_by, rem)
# left to fill the gap.
$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
